$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "schubert-winterreise_40"
$ws.Range("B2").Value = "schubert-winterreise_202"
$ws.Range("C2").Value = 0.4017857142857143
$ws.Range("D2").Value = "[['D:maj', 'A:7', 'D:maj'], ['D:maj/F#', 'G:maj', 'D:maj']]"
$ws.Range("E2").Value = "[['F:maj/C', 'C:7', 'F:maj'], ['F:maj/A', 'A#:maj', 'F:maj/C']]"
$ws.Range("F2").Value = "[(2.36, 11.66), (60.04, 67.08)]"
$ws.Range("G2").Value = "[(63.1, 64.32), (61.48, 63.58)]"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""

$ws.Range("A3").Value = "schubert-winterreise_61"
$ws.Range("B3").Value = "schubert-winterreise_48"
$ws.Range("C3").Value = 0.0945054945054945
$ws.Range("D3").Value = "[['G:maj', 'G:7/F', 'C:maj/E'], ['G:7/F', 'C:maj/E', 'G:maj/D'], ['G:maj', 'D:7/C', 'G:maj/B'], ['D:7', 'G:maj', 'G:maj']]"
$ws.Range("E3").Value = "[['F:maj', 'F:7', 'A#:maj'], ['F:7/D#', 'A#:maj/D', 'F:maj/C'], ['F:maj/C', 'C:7', 'F:maj'], ['C:7', 'F:maj', 'F:maj']]"
$ws.Range("F3").Value = "[(54.76, 62.82), (59.62, 64.34), (70.98, 80.44), (7.8, 13.6)]"
$ws.Range("G3").Value = "[(14.48, 22.82), (57.6, 60.72), (59.5, 65.04), (60.72, 68.12)]"
$ws.Range("H3").Value = "spotify:track:68YORkKP9uvlOQFMZZZwH5"
$ws.Range("I3").Value = ""

$ws.Range("A4").Value = "isophonics_147"
$ws.Range("B4").Value = "isophonics_76"
$ws.Range("C4").Value = 0.1418067226890756
$ws.Range("D4").Value = "[['A', 'E', 'B']]"
$ws.Range("E4").Value = "[['C', 'G', 'D/3']]"
$ws.Range("F4").Value = "[(23.861, 34.866)]"
$ws.Range("G4").Value = "[(9.336, 13.342)]"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""

$ws.Range("A5").Value = "isophonics_135"
$ws.Range("B5").Value = "isophonics_231"
$ws.Range("C5").Value = 0.09064112011790715
$ws.Range("D5").Value = "[['B', 'E', 'A'], ['A', 'B', 'E']]"
$ws.Range("E5").Value = "[['D/5', 'G', 'C'], ['C', 'D', 'G']]"
$ws.Range("F5").Value = "[(13.8465, 18.745911), (13.393711, 17.886772)]"
$ws.Range("G5").Value = "[(38.435396, 43.311587), (17.014988, 20.892721)]"
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = "spotify:track:4F1AgKpuFRMLEgtPETVwZk"

$ws.Range("A6").Value = "schubert-winterreise_109"
$ws.Range("B6").Value = "schubert-winterreise_44"
$ws.Range("C6").Value = 0.09980620155038759
$ws.Range("D6").Value = "[['D:min', 'D:min', 'A:7', 'D:min', 'D:min']]"
$ws.Range("E6").Value = "[['A#:min', 'A#:min/F', 'F:7', 'A#:min', 'A#:min']]"
$ws.Range("F6").Value = "[(14.22, 42.42)]"
$ws.Range("G6").Value = "[(23.5, 31.24)]"
$ws.Range("H6").Value = "spotify:track:5UYEp9kllA47IhttiiMuJ0"
$ws.Range("I6").Value = ""

$ws.Range("A7").Value = "schubert-winterreise_41"
$ws.Range("B7").Value = "schubert-winterreise_203"
$ws.Range("C7").Value = 0.2125874125874126
$ws.Range("D7").Value = "[['B:min', 'F#:7/C#', 'B:min/D', 'B:maj/D#']]"
$ws.Range("E7").Value = "[['G:min', 'D:7', 'G:min', 'G:maj']]"
$ws.Range("F7").Value = "[(0.66, 3.28)]"
$ws.Range("G7").Value = "[(48.42, 58.56)]"
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = "spotify:track:68YORkKP9uvlOQFMZZZwH5"

$ws.Range("A8").Value = "isophonics_21"
$ws.Range("B8").Value = "isophonics_136"
$ws.Range("C8").Value = 0.1018027571580064
$ws.Range("D8").Value = "[['C', 'C/b7', 'F']]"
$ws.Range("E8").Value = "[['A', 'A', 'D/5']]"
$ws.Range("F8").Value = "[(35.532, 40.124)]"
$ws.Range("G8").Value = "[(7.448, 10.553)]"
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = ""

$ws.Range("A9").Value = "schubert-winterreise_183"
$ws.Range("B9").Value = "schubert-winterreise_128"
$ws.Range("C9").Value = 0.323076923076923
$ws.Range("D9").Value = "[['D:maj/A', 'G:maj', 'D:maj/A']]"
$ws.Range("E9").Value = "[['G:maj', 'C:maj/G', 'G:maj']]"
$ws.Range("F9").Value = "[(128.54, 132.84)]"
$ws.Range("G9").Value = "[(18.32, 25.82)]"
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = "spotify:track:68YORkKP9uvlOQFMZZZwH5"

$ws.Range("A10").Value = "schubert-winterreise_215"
$ws.Range("B10").Value = "schubert-winterreise_108"
$ws.Range("C10").Value = 0.1517241379310345
$ws.Range("D10").Value = "[['G:min', 'D:7/G', 'G:min']]"
$ws.Range("E10").Value = "[['A:min', 'E:7', 'A:min']]"
$ws.Range("F10").Value = "[(15.78, 21.28)]"
$ws.Range("G10").Value = "[(14.84, 22.2)]"
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = "spotify:track:3OD2uwEUQKg0WyW9Lewata"

$ws.Range("A11").Value = "schubert-winterreise_136"
$ws.Range("B11").Value = "jaah_30"
$ws.Range("C11").Value = 0.1214285714285714
$ws.Range("D11").Value = "[['F:maj/C', 'C:7', 'F:maj'], ['F:maj', 'C:maj', 'F:maj']]"
$ws.Range("E11").Value = "[['Bb', 'F:7', 'Bb'], ['Bb', 'F', 'Bb']]"
$ws.Range("F11").Value = "[(79.04, 86.54), (2.5, 26.2)]"
$ws.Range("G11").Value = "[(12.51, 13.91), (0.46, 3.22)]"
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""

$ws.Range("A12").Value = "isophonics_4"
$ws.Range("B12").Value = "isophonics_191"
$ws.Range("C12").Value = 0.1079545454545455
$ws.Range("D12").Value = "[['G:7', 'C:min', 'F:min']]"
$ws.Range("E12").Value = "[['Bb:7', 'Eb:min', 'Ab:min']]"
$ws.Range("F12").Value = "[(156.152, 159.022)]"
$ws.Range("G12").Value = "[(8.620975, 14.367913)]"
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""

$ws.Range("A13").Value = "isophonics_1"
$ws.Range("B13").Value = "isophonics_82"
$ws.Range("C13").Value = 0.07037643207855973
$ws.Range("D13").Value = "[['Eb', 'Ab/5', 'Eb'], ['Bb', 'C', 'F']]"
$ws.Range("E13").Value = "[['A', 'D', 'A'], ['D', 'E', 'A']]"
$ws.Range("F13").Value = "[(17.016, 22.841), (40.03, 43.842)]"
$ws.Range("G13").Value = "[(3.988594, 7.754783), (52.750498, 59.809365)]"
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = "spotify:track:5EzvwjFwdP5Kfl5AZAemzu"

$ws.Range("A14").Value = "isophonics_132"
$ws.Range("B14").Value = "jaah_69"
$ws.Range("C14").Value = 0.07463144963144963
$ws.Range("D14").Value = "[['B', 'B', 'B/7']]"
$ws.Range("E14").Value = "[['Eb', 'Eb', 'Eb']]"
$ws.Range("F14").Value = "[(9.480113, 15.412811)]"
$ws.Range("G14").Value = "[(16.13, 24.98)]"
$ws.Range("H14").Value = ""
$ws.Range("I14").Value = ""

$ws.Range("A15").Value = "isophonics_166"
$ws.Range("B15").Value = "isophonics_295"
$ws.Range("C15").Value = 0.2708333333333333
$ws.Range("D15").Value = "[['D', 'G', 'G'], ['G', 'D', 'G']]"
$ws.Range("E15").Value = "[['G', 'C/5', 'C'], ['C/5', 'G', 'C/5']]"
$ws.Range("F15").Value = "[(19.139614, 27.777437), (17.839297, 22.785147)]"
$ws.Range("G15").Value = "[(20.870746, 26.837029), (12.094553, 16.853782)]"
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = ""

$ws.Range("A16").Value = "isophonics_180"
$ws.Range("B16").Value = "jaah_52"
$ws.Range("C16").Value = 0.1191151446398185
$ws.Range("D16").Value = "[['F', 'F:7', 'Bb', 'Bb:min', 'F']]"
$ws.Range("E16").Value = "[['F', 'F:7', 'Bb', 'Bb:min', 'F']]"
$ws.Range("F16").Value = "[(17.737518, 26.514661)]"
$ws.Range("G16").Value = "[(29.47, 34.95)]"
$ws.Range("H16").Value = ""
$ws.Range("I16").Value = ""

$ws.Range("A17").Value = "isophonics_99"
$ws.Range("B17").Value = "isophonics_241"
$ws.Range("C17").Value = 0.1164473684210526
$ws.Range("D17").Value = "[['G#', 'C#:min', 'F#:min'], ['F#:min', 'B', 'E']]"
$ws.Range("E17").Value = "[['G', 'C:min', 'F:min'], ['D:min', 'G', 'C']]"
$ws.Range("F17").Value = "[(2.804376, 6.298979), (16.376439, 18.860975)]"
$ws.Range("G17").Value = "[(38.783, 45.674), (90.971, 96.81)]"
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = ""
